$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row at row 100, pushing the existing rows
# (100:134) down to (101:135). This mirrors the diff, where the row
# that used to be at 100 is now at 101, ..., and the old row 134 is now
# row 135, while the brand new record occupies row 100.
$ws.Rows("100:100").Insert()

$newRow = 100

$ws.Cells.Item($newRow, 1).Value = 1
$ws.Cells.Item($newRow, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($newRow, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($newRow, 4).Value = 44466
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 15
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100108
$ws.Cells.Item($newRow, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value = 100108006
$ws.Cells.Item($newRow, 10).Value = "Plátano"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Verde"
$ws.Cells.Item($newRow, 13).Value = 120
$ws.Cells.Item($newRow, 14).Value = 18000
$ws.Cells.Item($newRow, 15).Value = 19000
$ws.Cells.Item($newRow, 16).Value = 18500
$ws.Cells.Item($newRow, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item($newRow, 18).Value = "Ecuador"
$ws.Cells.Item($newRow, 19).Value = 925
$ws.Cells.Item($newRow, 20).Value = 20
